$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new FWHM data run row for sg_rr_100_027 2023-12-08 17-44-55.csv
$r = 60

$ws.Cells.Item($r, 1).Value = "sg_rr_100_027 2023-12-08 17-44-55.csv"
$ws.Cells.Item($r, 2).Value = 0.01
$ws.Cells.Item($r, 3).Value = 1000
$ws.Cells.Item($r, 4).Value = 5001
$ws.Cells.Item($r, 5).Value = 1530
$ws.Cells.Item($r, 6).Value = 1570
$ws.Cells.Item($r, 7).Value = 0.5
$ws.Cells.Item($r, 8).Value = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = 0.98153846153846003
$ws.Cells.Item($r, 11).Value = 0.0043858818636388196
$ws.Cells.Item($r, 12).Value = "yes(although possible end peak not found)"
$ws.Cells.Item($r, 13).Value = 0.126322698170279
$ws.Cells.Item($r, 14).Value = 0.0032111503732971001
$ws.Cells.Item($r, 15).Value = "reduced approx fsr a bit, to see if this had any affect on fsr calculation as above, half the approx fsr was quite close to actual calculated fsr."

# Update view state to match the saved workbook: scrolled so row 44 is at
# top, with the active selection on A36.
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("A36").Select()
